# The author re-uploaded the workbook; the functional part of the change is
# that the "Expected Result" cell on the SampleServiceNew_1_0 sheet had its
# escaped literal "\n" (backslash + n) turned into a real line break inside
# the XML sample text, and column D was widened so the now-wrapped/longer
# text is readable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SampleServiceNew_1_0")

# D2 ("Expected Result" for the GET success-case row): replace the literal
# backslash-n with an actual newline between the XML declaration and the
# <SuccessResponse> element.
$ws.Range("D2").Value = "<?xml version=`"1.0`" encoding=`"UTF-8`"?>`n<SuccessResponse>MSISDN is proper</SuccessResponse>"

# Column D grew noticeably wider (from ~42 to ~62 characters) to accommodate
# the (now two-line) response text.
$ws.Columns.Item(4).ColumnWidth = 61.5
